$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 7576654
$ws.Range("I112").Value = 557.1429000000001
$ws.Range("J112").Value = 9405367
$ws.Range("K112").Value = 1671.4287
$ws.Range("L112").Value = 28216101
$ws.Range("M112").Value = -563.4287000000002
$ws.Range("N112").Value = -28218317
$ws.Range("H116").Value = 17305112
$ws.Range("I116").Value = 23072398
$ws.Range("J116").Value = 3250
$ws.Range("K116").Value = 23072398
$ws.Range("L116").Value = 3250
$ws.Range("M116").Value = -23068956
$ws.Range("N116").Value = -10134
$ws.Range("H121").Value = 1105
$ws.Range("J121").Value = 1105
$ws.Range("L121").Value = 3315
$ws.Range("N121").Value = -6809
$ws.Range("H132").Value = 195120.28
$ws.Range("I132").Value = 206469.36
$ws.Range("J132").Value = 61201.2
$ws.Range("K132").Value = 619408.08
$ws.Range("L132").Value = 183603.6
$ws.Range("M132").Value = -616878.08
$ws.Range("N132").Value = -188663.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 45381
$ws.Range("I2").Value = 60221.35
$ws.Range("J2").Value = 3333.3333
$ws.Range("K2").Value = 60221.35
$ws.Range("L2").Value = 3333.3333
$ws.Range("M2").Value = -60108.35
$ws.Range("N2").Value = -3559.3333
$ws.Range("H32").Value = 21373.871
$ws.Range("I32").Value = 2849.5193
$ws.Range("K32").Value = 2849.5193
$ws.Range("M32").Value = -2562.5193
$ws.Range("H61").Value = 1878.7291
$ws.Range("I61").Value = 1297.0264
$ws.Range("J61").Value = 4089.2
$ws.Range("K61").Value = 1297.0264
$ws.Range("L61").Value = 4089.2
$ws.Range("M61").Value = -1085.0264
$ws.Range("N61").Value = -4513.2
$ws.Range("H74").Value = 4997.2573
$ws.Range("I74").Value = 1333.4615
$ws.Range("J74").Value = 15581.556
$ws.Range("K74").Value = 1333.4615
$ws.Range("L74").Value = 15581.556
$ws.Range("M74").Value = -459.4614999999999
$ws.Range("N74").Value = -17329.556
$ws.Range("H77").Value = 4997.2573
$ws.Range("I77").Value = 1333.4615
$ws.Range("J77").Value = 15581.556
$ws.Range("K77").Value = 6667.307499999999
$ws.Range("L77").Value = 77907.78
$ws.Range("M77").Value = -2299.307499999999
$ws.Range("N77").Value = -86643.78
$ws.Range("H116").Value = 45381
$ws.Range("I116").Value = 60221.35
$ws.Range("J116").Value = 3333.3333
$ws.Range("K116").Value = 60221.35
$ws.Range("L116").Value = 3333.3333
$ws.Range("M116").Value = -57927.35
$ws.Range("N116").Value = -7921.3333
$ws.Range("H122").Value = 2496.2
$ws.Range("I122").Value = 2406.4707
$ws.Range("K122").Value = 7219.4121
$ws.Range("M122").Value = -4769.4121
$ws.Range("H132").Value = 2670.1316
$ws.Range("I132").Value = 2096.6453
$ws.Range("J132").Value = 5209.857
$ws.Range("K132").Value = 6289.9359
$ws.Range("L132").Value = 15629.571
$ws.Range("M132").Value = -3759.9359
$ws.Range("N132").Value = -20689.571
$ws.Range("H136").Value = 1878.7291
$ws.Range("I136").Value = 1297.0264
$ws.Range("J136").Value = 4089.2
$ws.Range("K136").Value = 3891.0792
$ws.Range("L136").Value = 12267.6
$ws.Range("M136").Value = -1341.0792
$ws.Range("N136").Value = -17367.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 45381
$ws.Range("I3").Value = 60221.35
$ws.Range("J3").Value = 3333.3333
$ws.Range("K3").Value = 60221.35
$ws.Range("L3").Value = 3333.3333
$ws.Range("M3").Value = -60107.35
$ws.Range("N3").Value = -3561.3333
$ws.Range("H20").Value = 4000
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -4494
$ws.Range("H99").Value = 1724.6786
$ws.Range("I99").Value = 1494.7826
$ws.Range("J99").Value = 2782.2
$ws.Range("K99").Value = 1494.7826
$ws.Range("L99").Value = 2782.2
$ws.Range("M99").Value = 3.217399999999998
$ws.Range("N99").Value = -5778.2
$ws.Range("H134").Value = 3622.0857
$ws.Range("I134").Value = 2579.0952
$ws.Range("J134").Value = 5186.5713
$ws.Range("K134").Value = 7737.285600000001
$ws.Range("L134").Value = 15559.7139
$ws.Range("M134").Value = -5202.285600000001
$ws.Range("N134").Value = -20629.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2261.25
$ws.Range("I122").Value = 1229.4546
$ws.Range("J122").Value = 3522.3333
$ws.Range("K122").Value = 3688.3638
$ws.Range("L122").Value = 10566.9999
$ws.Range("M122").Value = -1238.3638
$ws.Range("N122").Value = -15466.9999
$ws.Range("H132").Value = 2243.2068
$ws.Range("I132").Value = 1794.8
$ws.Range("J132").Value = 3795.3845
$ws.Range("K132").Value = 5384.4
$ws.Range("L132").Value = 11386.1535
$ws.Range("M132").Value = -2854.4
$ws.Range("N132").Value = -16446.1535
$ws.Range("H134").Value = 2511.5625
$ws.Range("I134").Value = 1473.2703
$ws.Range("J134").Value = 6004
$ws.Range("K134").Value = 4419.810899999999
$ws.Range("L134").Value = 18012
$ws.Range("M134").Value = -1884.810899999999
$ws.Range("N134").Value = -23082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1220.8387
$ws.Range("I5").Value = 551.5
$ws.Range("K5").Value = 1654.5
$ws.Range("M5").Value = -1542.5
$ws.Range("H135").Value = 1220.8387
$ws.Range("I135").Value = 551.5
$ws.Range("K135").Value = 4963.5
$ws.Range("M135").Value = -2428.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1390246
$ws.Range("I122").Value = 1390246
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4170738
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4168288
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2020.7073
$ws.Range("I132").Value = 1783.9546
$ws.Range("J132").Value = 2997.3125
$ws.Range("K132").Value = 5351.8638
$ws.Range("L132").Value = 8991.9375
$ws.Range("M132").Value = -2821.8638
$ws.Range("N132").Value = -14051.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1394.0714
$ws.Range("I16").Value = 1702.3
$ws.Range("K16").Value = 1702.3
$ws.Range("M16").Value = -1532.3
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H99").Value = 29500
$ws.Range("I99").Value = 22000
$ws.Range("K99").Value = 22000
$ws.Range("M99").Value = -19005
$ws.Range("H132").Value = 3885.2144
$ws.Range("I132").Value = 3744.0908
$ws.Range("J132").Value = 3976.5293
$ws.Range("K132").Value = 11232.2724
$ws.Range("L132").Value = 11929.5879
$ws.Range("M132").Value = -8702.2724
$ws.Range("N132").Value = -16989.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1463.2
$ws.Range("J23").Value = 2880.5
$ws.Range("L23").Value = 2880.5
$ws.Range("N23").Value = -3338.5
$ws.Range("H61").Value = 1686.6666
$ws.Range("I61").Value = 1686.6666
$ws.Range("K61").Value = 1686.6666
$ws.Range("M61").Value = -1394.6666
$ws.Range("H122").Value = 41741.24
$ws.Range("I122").Value = 67928.87
$ws.Range("J122").Value = 2459.8
$ws.Range("K122").Value = 203786.61
$ws.Range("L122").Value = 7379.400000000001
$ws.Range("M122").Value = -201336.61
$ws.Range("N122").Value = -12279.4
$ws.Range("H132").Value = 10871957
$ws.Range("I132").Value = 16668842
$ws.Range("J132").Value = 2797.3125
$ws.Range("K132").Value = 50006526
$ws.Range("L132").Value = 8391.9375
$ws.Range("M132").Value = -50003996
$ws.Range("N132").Value = -13451.9375
$ws.Range("H136").Value = 7114673.5
$ws.Range("I136").Value = 9288058
$ws.Range("J136").Value = 1778.6364
$ws.Range("K136").Value = 27864174
$ws.Range("L136").Value = 5335.9092
$ws.Range("M136").Value = -27861624
$ws.Range("N136").Value = -10435.9092
